$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 1).Value = './model_output/2025-08-22-14-58-15-Onefitall_11'
$ws.Cells.Item(12, 2).Value = 0.6342857142857141
$ws.Cells.Item(13, 1).Value = './model_output/2025-08-22-14-58-32-Onefitall_12'
$ws.Cells.Item(13, 2).Value = 0.6364285714285715
$ws.Cells.Item(14, 1).Value = './model_output/2025-08-22-14-58-48-Onefitall_13'
$ws.Cells.Item(14, 2).Value = 0.6721428571428572
$ws.Cells.Item(15, 1).Value = './model_output/2025-08-22-14-59-01-LLMFlareNet_1'
$ws.Cells.Item(15, 2).Value = 0.4842857142857143
$ws.Cells.Item(16, 1).Value = './model_output/2025-08-22-14-59-51-LLMFlareNet_2'
$ws.Cells.Item(16, 2).Value = 0.6464285714285715
$ws.Cells.Item(17, 1).Value = './model_output/2025-08-22-15-05-56-Onefitall_11'
$ws.Cells.Item(17, 2).Value = 0.5428571428571429
$ws.Cells.Item(18, 1).Value = './model_output/2025-08-22-15-06-00-Onefitall_12'
$ws.Cells.Item(18, 2).Value = 0.5357142857142858
$ws.Cells.Item(19, 1).Value = './model_output/2025-08-22-15-06-04-Onefitall_13'
$ws.Cells.Item(19, 2).Value = 0.6285714285714286
$ws.Cells.Item(20, 1).Value = './model_output/2025-08-22-15-06-07-LLMFlareNet_1'
$ws.Cells.Item(20, 2).Value = 0.4500000000000001
$ws.Cells.Item(21, 1).Value = './model_output/2025-08-22-15-06-14-LLMFlareNet_2'
$ws.Cells.Item(21, 2).Value = 0.5571428571428572
$ws.Cells.Item(22, 1).Value = './model_output/2025-08-22-15-07-56-Onefitall_11'
$ws.Cells.Item(22, 2).Value = 0.5428571428571429
$ws.Cells.Item(23, 1).Value = './model_output/2025-08-22-15-08-00-Onefitall_12'
$ws.Cells.Item(23, 2).Value = 0.5357142857142858
$ws.Cells.Item(24, 1).Value = './model_output/2025-08-22-15-08-04-Onefitall_13'
$ws.Cells.Item(24, 2).Value = 0.6285714285714286
$ws.Cells.Item(25, 1).Value = './model_output/2025-08-22-15-08-07-LLMFlareNet_1'
$ws.Cells.Item(25, 2).Value = 0.4500000000000001
$ws.Cells.Item(26, 1).Value = './model_output/2025-08-22-15-08-15-LLMFlareNet_2'
$ws.Cells.Item(26, 2).Value = 0.5571428571428572
$ws.Cells.Item(27, 1).Value = './model_output/2025-08-22-15-09-28-Onefitall_11'
$ws.Cells.Item(27, 2).Value = 0.5428571428571429
$ws.Cells.Item(28, 1).Value = './model_output/2025-08-22-15-09-32-Onefitall_12'
$ws.Cells.Item(28, 2).Value = 0.5357142857142858
$ws.Cells.Item(29, 1).Value = './model_output/2025-08-22-15-09-36-Onefitall_13'
$ws.Cells.Item(29, 2).Value = 0.6285714285714286
$ws.Cells.Item(30, 1).Value = './model_output/2025-08-22-15-09-40-LLMFlareNet_1'
$ws.Cells.Item(30, 2).Value = 0.4500000000000001
$ws.Cells.Item(31, 1).Value = './model_output/2025-08-22-15-09-47-LLMFlareNet_2'
$ws.Cells.Item(31, 2).Value = 0.5571428571428572
$ws.Cells.Item(32, 1).Value = './model_output/2025-08-22-15-11-24-Onefitall_11'
$ws.Cells.Item(32, 2).Value = 0.5428571428571429
$ws.Cells.Item(33, 1).Value = './model_output/2025-08-22-15-11-28-Onefitall_12'
$ws.Cells.Item(33, 2).Value = 0.5357142857142858
$ws.Cells.Item(34, 1).Value = './model_output/2025-08-22-15-11-32-Onefitall_13'
$ws.Cells.Item(34, 2).Value = 0.6285714285714286
$ws.Cells.Item(35, 1).Value = './model_output/2025-08-22-15-11-36-LLMFlareNet_1'
$ws.Cells.Item(35, 2).Value = 0.4500000000000001
$ws.Cells.Item(36, 1).Value = './model_output/2025-08-22-15-11-43-LLMFlareNet_2'
$ws.Cells.Item(36, 2).Value = 0.5571428571428572
$ws.Cells.Item(37, 1).Value = './model_output/2025-08-22-15-28-21-Onefitall_11'
$ws.Cells.Item(37, 2).Value = 0.5428571428571429
$ws.Cells.Item(38, 1).Value = './model_output/2025-08-22-15-28-25-Onefitall_12'
$ws.Cells.Item(38, 2).Value = 0.5357142857142858
$ws.Cells.Item(39, 1).Value = './model_output/2025-08-22-15-28-29-Onefitall_13'
$ws.Cells.Item(39, 2).Value = 0.6285714285714286
$ws.Cells.Item(40, 1).Value = './model_output/2025-08-22-15-28-47-Onefitall_11'
$ws.Cells.Item(40, 2).Value = 0.5428571428571429
$ws.Cells.Item(41, 1).Value = './model_output/2025-08-22-15-28-51-Onefitall_12'
$ws.Cells.Item(41, 2).Value = 0.5357142857142858
$ws.Cells.Item(42, 1).Value = './model_output/2025-08-22-15-28-55-Onefitall_13'
$ws.Cells.Item(42, 2).Value = 0.6285714285714286
$ws.Cells.Item(43, 1).Value = './model_output/2025-08-22-15-28-59-LLMFlareNet_1'
$ws.Cells.Item(43, 2).Value = 0.4500000000000001
$ws.Cells.Item(44, 1).Value = './model_output/2025-08-22-15-29-06-LLMFlareNet_2'
$ws.Cells.Item(44, 2).Value = 0.5571428571428572
$ws.Cells.Item(45, 1).Value = './model_output/2025-08-22-15-29-10-Onefitall_11'
$ws.Cells.Item(45, 2).Value = 0.6071428571428572
$ws.Cells.Item(46, 1).Value = './model_output/2025-08-22-15-29-14-Onefitall_12'
$ws.Cells.Item(46, 2).Value = 0.6071428571428572
$ws.Cells.Item(47, 1).Value = './model_output/2025-08-22-15-29-18-Onefitall_13'
$ws.Cells.Item(47, 2).Value = 0.6499999999999999
$ws.Cells.Item(48, 1).Value = './model_output/2025-08-22-15-29-22-LLMFlareNet_1'
$ws.Cells.Item(48, 2).Value = 0.4571428571428572
$ws.Cells.Item(49, 1).Value = './model_output/2025-08-22-15-29-29-LLMFlareNet_2'
$ws.Cells.Item(49, 2).Value = 0.5642857142857143

$wb.Save()
